# Conserto do erro com o rotulo da coluna 2050 nas tabelas e retirada
# das linhas com total das tabelas.
#
# Sheets 1-4: fix the mislabeled E1 header (was a stray numeric value,
# should read the next period label "2050" / "2041-2050") and remove the
# trailing "Total" row (row 13).
# Sheet 5: only needs the E1 header label fixed (it never had a Total row).
# Sheet 6: only needs its trailing "Total" row (row 4) removed.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Potencia Acumulada - SIN (MW)" ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "2050"
$ws.Rows.Item(13).Delete()

# --- Sheet 2: "Geracao Periodo Medio (MWMed)" ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "2050"
$ws.Rows.Item(13).Delete()

# --- Sheet 3: "Atendimento a Ponta(MW)" ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "2050"
$ws.Rows.Item(13).Delete()

# --- Sheet 4: "Potencia Incremental - SIN(MW)" (label is a range "2041-2050") ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "2041-2050"
$ws.Rows.Item(13).Delete()

# --- Sheet 5: "Emissoes Totais (MtCO2eq)" (no Total row to remove) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "2050"

# --- Sheet 6: "Custo Total (bilhoes de R$)" (only the Total row goes away) ---
$ws = $wb.Worksheets.Item(6)
$ws.Rows.Item(4).Delete()
